$wb = $excel.ActiveWorkbook

# Delete the extra sheets "phpmyadmin1" and "phpmyadmin2"
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("phpmyadmin2").Delete()
$wb.Worksheets.Item("phpmyadmin1").Delete()

# Update labels on the remaining sheet
$ws = $wb.Worksheets.Item("phpmyadmin")
$ws.Range("A2").Value = "V-php:S2068"
$ws.Range("A3").Value = "V-javascript:S2819"
$ws.Range("A4").Value = "V-php:S2964"
